$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 3.479168461094797
$ws.Range("D2").Value = 3.917626312918256
$ws.Range("E2").Value = 16.58652060621377
$ws.Range("F2").Value = 18.09493953642914
$ws.Range("G2").Value = 18.10443107543747
$ws.Range("H2").Value = 11.61406522871387
$ws.Range("I2").Value = 15.58329534938068
$ws.Range("K2").Value = 15.28291200518011
$ws.Range("O2").Value = 16.24382269447874
$ws.Range("C3").Value = 3.334030883129838
$ws.Range("D3").Value = 3.824017570080299
$ws.Range("E3").Value = 15.64034333264192
$ws.Range("F3").Value = 18.17882040490029
$ws.Range("G3").Value = 18.23446903175893
$ws.Range("H3").Value = 11.6869573045741
$ws.Range("I3").Value = 15.67646889796573
$ws.Range("K3").Value = 14.44354404531966
$ws.Range("O3").Value = 16.36920664856004
$ws.Range("C4").Value = 3.241884465073563
$ws.Range("D4").Value = 3.764962751397979
$ws.Range("E4").Value = 15.03397231761195
$ws.Range("F4").Value = 18.238730847159
$ws.Range("G4").Value = 18.3275543380983
$ws.Range("H4").Value = 11.73459581261667
$ws.Range("I4").Value = 15.73959680174709
$ws.Range("K4").Value = 13.90105350252084
$ws.Range("O4").Value = 16.4522880202233
$ws.Range("C5").Value = 3.203624529060056
$ws.Range("D5").Value = 3.740524086403425
$ws.Range("E5").Value = 14.78074727210228
$ws.Range("F5").Value = 18.26524059356407
$ws.Range("G5").Value = 18.36876672000908
$ws.Range("H5").Value = 11.75473273051641
$ws.Range("I5").Value = 15.76680026855475
$ws.Range("K5").Value = 13.6733117896379
$ws.Range("O5").Value = 16.48766811150018
$ws.Range("C6").Value = 3.197230315997795
$ws.Range("D6").Value = 3.736444234958253
$ws.Range("E6").Value = 14.73833841836699
$ws.Range("F6").Value = 18.26976854470446
$ws.Range("G6").Value = 18.37580657105059
$ws.Range("H6").Value = 11.75812014080846
$ws.Range("I6").Value = 15.77140637631168
$ws.Range("K6").Value = 13.63509733905048
$ws.Range("O6").Value = 16.49363472261561
$ws.Range("C7").Value = 3.241371276384765
$ws.Range("D7").Value = 3.764634642393869
$ws.Range("E7").Value = 15.0305816402524
$ws.Range("F7").Value = 18.23907990607315
$ws.Range("G7").Value = 18.32809693233771
$ws.Range("H7").Value = 11.73486445662579
$ws.Range("I7").Value = 15.73995770486408
$ws.Range("K7").Value = 13.89800889882719
$ws.Range("O7").Value = 16.45275901058077
$ws.Range("C8").Value = 3.429783779463988
$ws.Range("D8").Value = 3.885691378627899
$ws.Range("E8").Value = 16.26569776804289
$ws.Range("F8").Value = 18.12210596255715
$ws.Range("G8").Value = 18.14649012048138
$ws.Range("H8").Value = 11.63859955536826
$ws.Range("I8").Value = 15.61418709479035
$ws.Range("K8").Value = 14.99920400270334
$ws.Range("O8").Value = 16.28578463901969
$ws.Range("C9").Value = 3.773216698685664
$ws.Range("D9").Value = 4.109596906727188
$ws.Range("E9").Value = 18.59163789771513
$ws.Range("F9").Value = 17.9602064399366
$ws.Range("G9").Value = 17.89758358237699
$ws.Range("H9").Value = 11.47274059897726
$ws.Range("I9").Value = 15.41496040492272
$ws.Range("K9").Value = 16.93850863735358
$ws.Range("O9").Value = 16.0071153267247
$ws.Range("C10").Value = 4.007407071021973
$ws.Range("D10").Value = 4.264702393678498
$ws.Range("E10").Value = 20.2435151697773
$ws.Range("F10").Value = 17.88339917965543
$ws.Range("G10").Value = 17.78287062533785
$ws.Range("H10").Value = 11.36491790857209
$ws.Range("I10").Value = 15.29809172234518
$ws.Range("K10").Value = 18.22405192039199
$ws.Range("O10").Value = 15.83267095950727
$ws.Range("C11").Value = 4.11392849930784
$ws.Range("D11").Value = 4.333005233089537
$ws.Range("E11").Value = 20.95250246555571
$ws.Range("F11").Value = 17.85779950714607
$ws.Range("G11").Value = 17.74602708329904
$ws.Range("H11").Value = 11.31892895191193
$ws.Range("I11").Value = 15.2514576902609
$ws.Range("K11").Value = 18.77793647231413
$ws.Range("O11").Value = 15.76001468908334
$ws.Range("C12").Value = 4.155743757859213
$ws.Range("D12").Value = 4.358529204140212
$ws.Range("E12").Value = 21.21490853597217
$ws.Range("F12").Value = 17.84946261183281
$ws.Range("G12").Value = 17.73432122259554
$ws.Range("H12").Value = 11.30195563854382
$ws.Range("I12").Value = 15.23474779383734
$ws.Range("K12").Value = 18.98319136502883
$ws.Range("O12").Value = 15.73347569991225
$ws.Range("C13").Value = 4.146777896284066
$ws.Range("D13").Value = 4.353047566226368
$ws.Range("E13").Value = 21.15866421536343
$ws.Range("F13").Value = 17.85119754300649
$ws.Range("G13").Value = 17.73674179835134
$ws.Range("H13").Value = 11.30559147291057
$ws.Range("I13").Value = 15.23830418119174
$ws.Range("K13").Value = 18.93918628516183
$ws.Range("O13").Value = 15.7391478338092
$ws.Range("C14").Value = 4.117386589821782
$ws.Range("D14").Value = 4.335112025694587
$ws.Range("E14").Value = 20.97421227147661
$ws.Range("F14").Value = 17.8570863509883
$ws.Range("G14").Value = 17.74501880531143
$ws.Range("H14").Value = 11.31752368494023
$ws.Range("I14").Value = 15.25006387094789
$ws.Range("K14").Value = 18.79491309572161
$ws.Range("O14").Value = 15.75781171909384
$ws.Range("C15").Value = 4.09926715823089
$ws.Range("D15").Value = 4.324081137898547
$ws.Range("E15").Value = 20.86044036717169
$ws.Range("F15").Value = 17.86087054960368
$ws.Range("G15").Value = 17.75038233377145
$ws.Range("H15").Value = 11.32489008328908
$ws.Range("I15").Value = 15.25739096475798
$ws.Range("K15").Value = 18.70595587638921
$ws.Range("O15").Value = 15.76937110836276
$ws.Range("C16").Value = 4.000637166813752
$ws.Range("D16").Value = 4.260191813313216
$ws.Range("E16").Value = 20.19632966225459
$ws.Range("F16").Value = 17.88526135804374
$ws.Range("G16").Value = 17.78559062830045
$ws.Range("H16").Value = 11.36798511305388
$ws.Range("I16").Value = 15.30127166027683
$ws.Range("K16").Value = 18.18722818067449
$ws.Range("O16").Value = 15.83755496226383
$ws.Range("C17").Value = 3.940822754900111
$ws.Range("D17").Value = 4.220407904388013
$ws.Range("E17").Value = 19.7780682639319
$ws.Range("F17").Value = 17.90262699979846
$ws.Range("G17").Value = 17.81114838050631
$ws.Range("H17").Value = 11.39520744881882
$ws.Range("I17").Value = 15.32987101124024
$ws.Range("K17").Value = 17.86105267000377
$ws.Range("O17").Value = 15.88110718418184
$ws.Range("C18").Value = 3.906015300692447
$ws.Range("D18").Value = 4.197313852843833
$ws.Range("E18").Value = 19.53349684052681
$ws.Range("F18").Value = 17.91349365609069
$ws.Range("G18").Value = 17.82729020703168
$ws.Range("H18").Value = 11.41115294540311
$ws.Range("I18").Value = 15.34693463810678
$ws.Range("K18").Value = 17.67053702900286
$ws.Range("O18").Value = 15.90678717448599
$ws.Range("C19").Value = 3.894161584912663
$ws.Range("D19").Value = 4.189458822683044
$ws.Range("E19").Value = 19.45000137802718
$ws.Range("F19").Value = 17.91732333653463
$ws.Range("G19").Value = 17.83300188185563
$ws.Range("H19").Value = 11.41660121927718
$ws.Range("I19").Value = 15.35281723097731
$ws.Range("K19").Value = 17.60553374975147
$ws.Range("O19").Value = 15.91558983233579
$ws.Range("C20").Value = 3.947232088439949
$ws.Range("D20").Value = 4.224664972626229
$ws.Range("E20").Value = 19.82300649213283
$ws.Range("F20").Value = 17.9006873808811
$ws.Range("G20").Value = 17.80827820317952
$ws.Range("H20").Value = 11.39227977329284
$ws.Range("I20").Value = 15.32676293517373
$ws.Range("K20").Value = 17.89607593981938
$ws.Range("O20").Value = 15.87640570389503
$ws.Range("C21").Value = 4.126043803092162
$ws.Range("D21").Value = 4.340389503112812
$ws.Range("E21").Value = 21.02855480466343
$ws.Range("F21").Value = 17.85531972769803
$ws.Range("G21").Value = 17.74252639059977
$ws.Range("H21").Value = 11.31400690358553
$ws.Range("I21").Value = 15.24658391690631
$ws.Range("K21").Value = 18.83741175577979
$ws.Range("O21").Value = 15.75230315100799
$ws.Range("C22").Value = 4.246085528761705
$ws.Range("D22").Value = 4.41402890263292
$ws.Range("E22").Value = 21.78107332677352
$ws.Range("F22").Value = 17.83358536947
$ws.Range("G22").Value = 17.71266107680089
$ws.Range("H22").Value = 11.26542649009976
$ws.Range("I22").Value = 15.19972078184904
$ws.Range("K22").Value = 19.42645001719906
$ws.Range("O22").Value = 15.6768795280164
$ws.Range("C23").Value = 4.182495709525043
$ws.Range("D23").Value = 4.374913633149544
$ws.Range("E23").Value = 21.38266596490005
$ws.Range("F23").Value = 17.84445680569581
$ws.Range("G23").Value = 17.72738914519061
$ws.Range("H23").Value = 11.29111854857986
$ws.Range("I23").Value = 15.22422237363608
$ws.Range("K23").Value = 19.11447622493725
$ws.Range("O23").Value = 15.71661071408861
$ws.Range("C24").Value = 3.94433573299126
$ws.Range("D24").Value = 4.222741041504537
$ws.Range("E24").Value = 19.80270271593401
$ws.Range("F24").Value = 17.90156153337551
$ws.Range("G24").Value = 17.80957130188604
$ws.Range("H24").Value = 11.39360245652267
$ws.Range("I24").Value = 15.32816616121954
$ws.Range("K24").Value = 17.88025123866248
$ws.Range("O24").Value = 15.87852924645903
$ws.Range("C25").Value = 3.683350909743294
$ws.Range("D25").Value = 4.05059579862901
$ws.Range("E25").Value = 17.94544277659102
$ws.Range("F25").Value = 17.99667256326437
$ws.Range("G25").Value = 17.95314272212814
$ws.Range("H25").Value = 11.51515020872947
$ws.Range("I25").Value = 15.46371848264007
$ws.Range("K25").Value = 16.43804160077563
$ws.Range("O25").Value = 16.07722437957745
